$wb = $excel.ActiveWorkbook

# The localization status report is being regenerated: items that were
# previously "Ready for handoff" are now shown as "In Translation".
$newStatus = "In Translation"

# Overview sheet: per-language status columns (zh-cn = E, de-de = F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus

# Per-language detail sheets: "Status" column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus

# The shorter status text lets the affected columns shrink, matching the
# width the report generator recalculates after refreshing the values.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
